# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" values (column G) for data rows 2-30, replacing the previous
# Strike# derived values with the freshly computed K values.
$kValues = @{
    2  = 4
    3  = 4
    4  = 3
    5  = 3
    6  = 3
    7  = 0
    8  = 3
    9  = 8
    10 = 9
    11 = 3
    12 = 8
    13 = 3
    14 = 6
    15 = 5
    16 = 4
    17 = 2
    18 = 6
    19 = 4
    20 = 5
    21 = 8
    22 = 4
    23 = 8
    24 = 3
    25 = 3
    26 = 4
    27 = 3
    28 = 4
    29 = 3
    30 = 0
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
